$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 2.75
$ws.Range("D3").Value = 2.5
$ws.Range("D4").Value = 1
$ws.Range("D6").Value = 1
